# Updates cryptos list prices (column D) and 1h volume % changes (column E)
# for rows 2-51, matching the latest data pull from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$origStyle = $c.Style
$c.Value = "'57.983.18"
$c.Style = $origStyle
$ws.Range("E2").Value = '  -0.79%  '

$c = $ws.Range("D3")
$origStyle = $c.Style
$c.Value = "'2.446.44"
$c.Style = $origStyle
$ws.Range("E3").Value = '  -3.05%  '

$ws.Range("E4").Value = '  +0.13%  '

$c = $ws.Range("D5")
$origStyle = $c.Style
$c.Value = "'525.09"
$c.Style = $origStyle
$ws.Range("E5").Value = '  +0.54%  '

$c = $ws.Range("D6")
$origStyle = $c.Style
$c.Value = "'131.10"
$c.Style = $origStyle
$ws.Range("E6").Value = '  -1.47%  '

$c = $ws.Range("D7")
$origStyle = $c.Style
$c.Value = "'0.999"
$c.Style = $origStyle
$ws.Range("E7").Value = '  -0.05%  '

$c = $ws.Range("D8")
$origStyle = $c.Style
$c.Value = "'0.564"
$c.Style = $origStyle
$ws.Range("E8").Value = '  -0.21%  '

$c = $ws.Range("D9")
$origStyle = $c.Style
$c.Value = "'2.450.63"
$c.Style = $origStyle
$ws.Range("E9").Value = '  -2.88%  '

$ws.Range("E10").Value = '  -0.23%  '

$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("E12").Value = '  -4.06%  '

$c = $ws.Range("D13")
$origStyle = $c.Style
$c.Value = "'0.324"
$c.Style = $origStyle
$ws.Range("E13").Value = '  -2.13%  '

$c = $ws.Range("D14")
$origStyle = $c.Style
$c.Value = "'2.880.40"
$c.Style = $origStyle
$ws.Range("E14").Value = '  -3.02%  '

$c = $ws.Range("D15")
$origStyle = $c.Style
$c.Value = "'57.864.79"
$c.Style = $origStyle
$ws.Range("E15").Value = '  -0.91%  '

$c = $ws.Range("D16")
$origStyle = $c.Style
$c.Value = "'21.76"
$c.Style = $origStyle
$ws.Range("E16").Value = '  -1.52%  '

$ws.Range("E17").Value = '  -1.65%  '

$c = $ws.Range("D18")
$origStyle = $c.Style
$c.Value = "'2.448.74"
$c.Style = $origStyle
$ws.Range("E18").Value = '  -3.00%  '

$c = $ws.Range("D19")
$origStyle = $c.Style
$c.Value = "'10.37"
$c.Style = $origStyle
$ws.Range("E19").Value = '  -2.75%  '

$ws.Range("E20").Value = '  -1.34%  '

$c = $ws.Range("D21")
$origStyle = $c.Style
$c.Value = "'311.35"
$c.Style = $origStyle
$ws.Range("E21").Value = '  -3.19%  '

$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("E23").Value = '  -0.07%  '

$c = $ws.Range("D24")
$origStyle = $c.Style
$c.Value = "'65.00"
$c.Style = $origStyle
$ws.Range("E24").Value = '  +0.27%  '

$c = $ws.Range("D25")
$origStyle = $c.Style
$c.Value = "'0.403"
$c.Style = $origStyle
$ws.Range("E25").Value = '  -0.97%  '

$c = $ws.Range("D26")
$origStyle = $c.Style
$c.Value = "'2.577.39"
$c.Style = $origStyle
$ws.Range("E26").Value = '  -2.17%  '

$c = $ws.Range("D27")
$origStyle = $c.Style
$c.Value = "'1.00"
$c.Style = $origStyle
$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("E28").Value = '  -1.61%  '

$c = $ws.Range("D29")
$origStyle = $c.Style
$c.Value = "'7.25"
$c.Style = $origStyle
$ws.Range("E29").Value = '  -1.84%  '

$c = $ws.Range("D30")
$origStyle = $c.Style
$c.Value = "'174.03"
$c.Style = $origStyle
$ws.Range("E30").Value = '  +3.10%  '

$ws.Range("E31").Value = '  -1.91%  '

$ws.Range("E32").Value = '  -1.45%  '

$ws.Range("E33").Value = '  -1.40%  '

$ws.Range("E34").Value = '  -4.14%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("E36").Value = '  -0.20%  '

$c = $ws.Range("D37")
$origStyle = $c.Style
$c.Value = "'17.79"
$c.Style = $origStyle
$ws.Range("E37").Value = '  -2.04%  '

$ws.Range("E38").Value = '  -4.96%  '

$ws.Range("E39").Value = '  -3.34%  '

$c = $ws.Range("D40")
$origStyle = $c.Style
$c.Value = "'0.816"
$c.Style = $origStyle
$ws.Range("E40").Value = '  +5.98%  '

$c = $ws.Range("D41")
$origStyle = $c.Style
$c.Value = "'36.23"
$c.Style = $origStyle
$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("E42").Value = '  -2.33%  '

$c = $ws.Range("D43")
$origStyle = $c.Style
$c.Value = "'3.40"
$c.Style = $origStyle
$ws.Range("E43").Value = '  -1.55%  '

$c = $ws.Range("D44")
$origStyle = $c.Style
$c.Value = "'261.76"
$c.Style = $origStyle
$ws.Range("E44").Value = '  -5.13%  '

$c = $ws.Range("D45")
$origStyle = $c.Style
$c.Value = "'0.585"
$c.Style = $origStyle
$ws.Range("E45").Value = '  -2.16%  '

$c = $ws.Range("D46")
$origStyle = $c.Style
$c.Value = "'4.80"
$c.Style = $origStyle
$ws.Range("E46").Value = '  -3.76%  '

$ws.Range("E47").Value = '  +0.06%  '

$c = $ws.Range("D48")
$origStyle = $c.Style
$c.Value = "'121.71"
$c.Style = $origStyle
$ws.Range("E48").Value = '  -6.10%  '

$c = $ws.Range("D49")
$origStyle = $c.Style
$c.Value = "'0.0494"
$c.Style = $origStyle
$ws.Range("E49").Value = '  -0.97%  '

$ws.Range("E50").Value = '  -1.09%  '

$ws.Range("E51").Value = '  -4.02%  '
